$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.019102109518226
$ws.Range("D2").Value = 1.024626414001997
$ws.Range("E2").Value = 1.020284784967306
$ws.Range("I2").Value = 1.027638643975003
$ws.Range("J2").Value = 1.024306778841953
$ws.Range("K2").Value = 1.027454363206155
$ws.Range("L2").Value = 1.023125524464673
$ws.Range("N2").Value = 1.012186590035529
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.020229525089141
$ws.Range("D3").Value = 1.025469254231171
$ws.Range("E3").Value = 1.021244341889865
$ws.Range("I3").Value = 1.027858840210695
$ws.Range("J3").Value = 1.025069872959936
$ws.Range("K3").Value = 1.02810423740301
$ws.Range("L3").Value = 1.023890866851306
$ws.Range("N3").Value = 1.012438781994684
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.020958311289886
$ws.Range("D4").Value = 1.026013471865299
$ws.Range("E4").Value = 1.02186503089431
$ws.Range("I4").Value = 1.027999021287577
$ws.Range("J4").Value = 1.025562425080397
$ws.Range("K4").Value = 1.028522956605283
$ws.Range("L4").Value = 1.024385275011612
$ws.Range("N4").Value = 1.012601529762549
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.021264520241685
$ws.Range("D5").Value = 1.026241984146324
$ws.Range("E5").Value = 1.022125918980891
$ws.Range("I5").Value = 1.028057402484114
$ws.Range("J5").Value = 1.025769202150844
$ws.Range("K5").Value = 1.028698557038836
$ws.Range("L5").Value = 1.024592928358285
$ws.Range("N5").Value = 1.012669844353428
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.021315923978235
$ws.Range("D6").Value = 1.026280336131621
$ws.Range("E6").Value = 1.022169720341388
$ws.Range("I6").Value = 1.028067172632381
$ws.Range("J6").Value = 1.025803903829072
$ws.Range("K6").Value = 1.028728015962949
$ws.Range("L6").Value = 1.024627782794602
$ws.Range("N6").Value = 1.012681308545096
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.020962403546163
$ws.Range("D7").Value = 1.026016526345293
$ws.Range("E7").Value = 1.021868517088079
$ws.Range("I7").Value = 1.027999803544985
$ws.Range("J7").Value = 1.025565189191983
$ws.Range("K7").Value = 1.028525304672949
$ws.Range("L7").Value = 1.024388050454918
$ws.Range("N7").Value = 1.012602442996212
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.019483276215139
$ws.Range("D8").Value = 1.024911495253439
$ws.Range("E8").Value = 1.020609115101198
$ws.Range("I8").Value = 1.027713536590511
$ws.Range("J8").Value = 1.024564923105443
$ws.Range("K8").Value = 1.027674362520364
$ws.Range("L8").Value = 1.023384345417709
$ws.Range("N8").Value = 1.012271910106137
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.016871237349404
$ws.Range("D9").Value = 1.022955430374831
$ws.Range("E9").Value = 1.018388274543187
$ws.Range("I9").Value = 1.027191483642956
$ws.Range("J9").Value = 1.022792956422812
$ws.Range("K9").Value = 1.026161156772768
$ws.Range("L9").Value = 1.021609395075031
$ws.Range("N9").Value = 1.011686114059901
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.015126007925561
$ws.Range("D10").Value = 1.021645422496597
$ws.Range("E10").Value = 1.016906602288602
$ws.Range("I10").Value = 1.026831611154439
$ws.Range("J10").Value = 1.021605313418894
$ws.Range("K10").Value = 1.025143104676516
$ws.Range("L10").Value = 1.020421835728536
$ws.Range("N10").Value = 1.011293318474671
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.014369365908703
$ws.Range("D11").Value = 1.021076756583766
$ws.Range("E11").Value = 1.016264750864665
$ws.Range("I11").Value = 1.026672975108853
$ws.Range("J11").Value = 1.021089540765051
$ws.Range("K11").Value = 1.024700079286251
$ws.Range("L11").Value = 1.019906592006171
$ws.Range("N11").Value = 1.011122693996241
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.014088171326678
$ws.Range("D12").Value = 1.020865314385518
$ws.Range("E12").Value = 1.016026296214914
$ws.Range("I12").Value = 1.026613628535422
$ws.Range("J12").Value = 1.020897731320956
$ws.Range("K12").Value = 1.02453518859487
$ws.Range("L12").Value = 1.019715053158203
$ws.Range("N12").Value = 1.011059234873245
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.014148495068179
$ws.Range("D13").Value = 1.020910679162257
$ws.Range("E13").Value = 1.016077447480054
$ws.Range("I13").Value = 1.026626377678437
$ws.Range("J13").Value = 1.020938885452387
$ws.Range("K13").Value = 1.024570573209667
$ws.Range("L13").Value = 1.01975614588974
$ws.Range("N13").Value = 1.011072850767039
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.014346125226402
$ws.Range("D14").Value = 1.021059283081168
$ws.Range("E14").Value = 1.016245040999359
$ws.Range("I14").Value = 1.026668078113452
$ws.Range("J14").Value = 1.021073690403167
$ws.Range("K14").Value = 1.024686456132337
$ws.Range("L14").Value = 1.019890762499063
$ws.Range("N14").Value = 1.011117450111343
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.014467872544728
$ws.Range("D15").Value = 1.021150814442062
$ws.Range("E15").Value = 1.016348295249737
$ws.Range("I15").Value = 1.026693715195065
$ws.Range("J15").Value = 1.021156717876356
$ws.Range("K15").Value = 1.024757811505323
$ws.Range("L15").Value = 1.019973683754602
$ws.Range("N15").Value = 1.011144918419505
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.015176203668229
$ws.Range("D16").Value = 1.021683132944284
$ws.Range("E16").Value = 1.016949193942022
$ws.Range("I16").Value = 1.026842080119871
$ws.Range("J16").Value = 1.021639511551934
$ws.Range("K16").Value = 1.025172460374958
$ws.Range("L16").Value = 1.020456009147723
$ws.Range("N16").Value = 1.011304630831031
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.015620266421239
$ws.Range("D17").Value = 1.02201666076676
$ws.Range("E17").Value = 1.017326046980078
$ws.Range("I17").Value = 1.026934393447652
$ws.Range("J17").Value = 1.021941948964639
$ws.Range("K17").Value = 1.02543196874615
$ws.Range("L17").Value = 1.020758284803533
$ws.Range("N17").Value = 1.011404669080987
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.015879189549798
$ws.Range("D18").Value = 1.022211064551961
$ws.Range("E18").Value = 1.017545832252202
$ws.Range("I18").Value = 1.026987967239525
$ws.Range("J18").Value = 1.022118209434547
$ws.Range("K18").Value = 1.025583123150602
$ws.Range("L18").Value = 1.020934498455427
$ws.Range("N18").Value = 1.011462967480086
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.01596746024019
$ws.Range("D19").Value = 1.022277327934847
$ws.Range("E19").Value = 1.017620768839136
$ws.Range("I19").Value = 1.027006188552247
$ws.Range("J19").Value = 1.022178284934623
$ws.Range("K19").Value = 1.025634626877875
$ws.Range("L19").Value = 1.020994566050902
$ws.Range("N19").Value = 1.01148283687796
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.015572632133119
$ws.Range("D20").Value = 1.021980890597651
$ws.Range("E20").Value = 1.017285616985508
$ws.Range("I20").Value = 1.026924517133547
$ws.Range("J20").Value = 1.021909515415958
$ws.Range("K20").Value = 1.025404147930662
$ws.Range("L20").Value = 1.020725863692647
$ws.Range("N20").Value = 1.011393941326604
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.014287932062354
$ws.Range("D21").Value = 1.021015528851882
$ws.Range("E21").Value = 1.016195690063099
$ws.Range("I21").Value = 1.026655810026954
$ws.Range("J21").Value = 1.021034000020107
$ws.Range("K21").Value = 1.02465234065046
$ws.Range("L21").Value = 1.019851125530745
$ws.Range("N21").Value = 1.011104318970547
$ws.Range("B22").Value = 1.019999999999999
$ws.Range("C22").Value = 1.013479355771098
$ws.Range("D22").Value = 1.020407327162554
$ws.Range("E22").Value = 1.015510163746761
$ws.Range("I22").Value = 1.026484421113991
$ws.Range("J22").Value = 1.020482206189651
$ws.Range("K22").Value = 1.024177732309568
$ws.Range("L22").Value = 1.019300249309252
$ws.Range("N22").Value = 1.010921749688158
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.013908077089011
$ws.Range("D23").Value = 1.020729864118113
$ws.Range("E23").Value = 1.015873597791345
$ws.Range("I23").Value = 1.026575509131655
$ws.Range("J23").Value = 1.020774848188867
$ws.Range("K23").Value = 1.02442951307544
$ws.Range("L23").Value = 1.019592364173738
$ws.Range("N23").Value = 1.01101857796621
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.015594156298629
$ws.Range("D24").Value = 1.021997054021466
$ws.Range("E24").Value = 1.017303885642032
$ws.Range("I24").Value = 1.026928980652009
$ws.Range("J24").Value = 1.021924171191776
$ws.Range("K24").Value = 1.025416719615761
$ws.Range("L24").Value = 1.02074051370153
$ws.Range("N24").Value = 1.01139878889828
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.01754718773915
$ws.Range("D25").Value = 1.023462170638335
$ws.Range("E25").Value = 1.018962608983072
$ws.Range("I25").Value = 1.027328533575969
$ws.Range("J25").Value = 1.023252165907
$ws.Range("K25").Value = 1.026553985818026
$ws.Range("L25").Value = 1.022069010905691
$ws.Range("N25").Value = 1.011837954880379
